$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab17")

# Fix mis-encoded accented characters in the "Regional Economic Communities" note
# (PALOP = "Países Africanos de Língua Oficial Portuguesa"; MERCOSUR = "Mercado Común del Sur")
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# Minor data correction on row 67
$ws.Range("I67").Value = 75378.546306000004

# Updated figures for row 97
$ws.Range("C97").Value = 144193.31095099999
$ws.Range("D97").Value = 78305.983385
$ws.Range("E97").Value = 22463.926414000001
$ws.Range("F97").Value = 244963.22075000001
$ws.Range("G97").Value = 19017.378799999999
$ws.Range("H97").Value = 159093.021412
$ws.Range("I97").Value = 94732.637432999996
$ws.Range("J97").Value = 272843.03764499997

# Updated figures for row 98
$ws.Range("C98").Value = 40674.361312000001
$ws.Range("D98").Value = 59146.964883000001
$ws.Range("E98").Value = 112031.16651700001
$ws.Range("F98").Value = 211852.49271200001
$ws.Range("G98").Value = 39684.129846999997
$ws.Range("H98").Value = 203628.562393
$ws.Range("I98").Value = 84104.150620999993
$ws.Range("J98").Value = 327416.84286099998

$wb.Save()
